# [Improvement] On terminology : room -> bed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rooms")

# Rename the worksheet itself
$ws.Name = "beds"

# Update the header-row terminology from "rooms" to "beds"
$ws.Range("A1").Value = "all_beds"
$ws.Range("B1").Value = "new_beds"
$ws.Range("C1").Value = "old_beds"
$ws.Range("E1").Value = "new_beds_service"
$ws.Range("F1").Value = "old_beds_service"
$ws.Range("G1").Value = "beds_capacities"

# Make the renamed "beds" sheet the active tab
$ws.Activate()
